{"js": "// Insert a new paragraph right after the paragraph that ends with\n// \"This means roughly 0.17 percent of all transactions.\" and before the\n// first empty paragraph that follows it. The new paragraph uses the same\n// paragraph formatting (1.5 line spacing) and run formatting (Times New\n// Roman, 12pt) as the surrounding text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the anchor paragraph by its distinctive text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"This means roughly 0.17 percent of all transactions.\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert the new (initially empty) paragraph right after the anchor.\nconst newPara = anchor.insertParagraph(\"\", \"After\");\n\n// The sentence is written out in the same run chunks as the original edit.\nconst chunks = [\n  \"In this paper I used multiple supervised learning \",\n  \"algorithms\",\n  \", deep learning models and compared their ROC_AUC score, F1-Score, Precision and Accuracy on the real\",\n  \"-\",\n  \"world dataset.\"\n];\n\nfor (const chunk of chunks) {\n  newPara.insertText(chunk, \"End\");\n}\nawait context.sync();\n\n// Apply the paragraph + run formatting (Times New Roman, 12pt, 1.5-line\n// spacing) to match the rest of the document.\nnewPara.lineSpacing = 18;\nnewPara.font.name = \"Times New Roman\";\nnewPara.font.nameBidirectional = \"Times New Roman\";\nnewPara.font.size = 12;\nnewPara.font.sizeBidirectional = 12;\n\nconst fullRange = newPara.getRange();\nfullRange.font.name = \"Times New Roman\";\nfullRange.font.nameBidirectional = \"Times New Roman\";\nfullRange.font.size = 12;\nfullRange.font.sizeBidirectional = 12;\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph right after the paragraph that ends with\n# \"This means roughly 0.17 percent of all transactions.\" and before the\n# first empty paragraph that follows it. Word copies the paragraph mark\n# formatting (Times New Roman, 12pt, 1.5-line spacing) from the\n# surrounding text automatically when the new paragraph is created, so\n# there is no need to re-apply it explicitly.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"This means roughly 0.17 percent of all transactions.\"\n\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $anchor = $p\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph not found\"\n}\n\n$anchor.Range.InsertParagraphAfter()\n\n$newPara = $anchor.Next()\n$newPara.Range.Text = \"In this paper I used multiple supervised learning algorithms, deep learning models and compared their ROC_AUC score, F1-Score, Precision and Accuracy on the real-world dataset.\"\n"}
